$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add a new log entry row (row 46) ---
# First, copy row 45's current formatting (the "most recent entry" highlighted
# style) down to row 46, so the new last row inherits that look.
$ws.Range("A45:C45").Copy()
$ws.Range("A46:C46").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's values: date (2025-04-02), hours, and description.
$ws.Range("A46").Value = 45749
$ws.Range("B46").Value = 2
$ws.Range("C46").Value = "Gathered evaluation from classmates and updated final report"

# Row 45 is no longer the last entry, so restore its formatting back to the
# plain (non-highlighted) style used by the rows above it.
$ws.Range("B44:C44").Copy()
$ws.Range("B45:C45").PasteSpecial(-4122)
$excel.CutCopyMode = 0
